$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Henry E Jones"
$ws.Range("B4").Value = "000-000-0000"
$ws.Range("C4").Value = "hjones@belhaven.edu"
